$d = $word.ActiveDocument
Write-Output ("Exists before: " + $d.Bookmarks.Exists("_GoBack"))
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    Write-Output "deleted called"
} catch {
    Write-Output ("ERROR1: " + $_.Exception.Message)
}
Write-Output ("Exists after: " + $d.Bookmarks.Exists("_GoBack"))
